# Updated cryptos list values (Price / Volume(1h)) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.881.51'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.496.09'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '531.51'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').Value = '134.93'
$ws.Range('E6').Value = '  -0.85%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').Value = '0.566'
$ws.Range('E8').Value = '  +1.32%  '
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('E10').Value = '  -0.93%  '
$ws.Range('D11').Value = '5.38'
$ws.Range('E11').Value = '  +1.54%  '
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').Value = '2.936.46'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('D14').Value = '58.802.72'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('E15').Value = '  -1.55%  '
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = '2.504.79'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '11.02'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').Value = '323.23'
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').Value = '5.94'
$ws.Range('E22').Value = '  +1.84%  '
$ws.Range('D23').Value = '64.97'
$ws.Range('E23').Value = '  +2.79%  '
$ws.Range('D24').Value = '0.419'
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('D25').Value = '0.163'
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.50'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').Value = '0.0₃0758'
$ws.Range('E28').Value = '  -1.39%  '
$ws.Range('D29').Value = '170.05'
$ws.Range('E29').Value = '  +1.62%  '
$ws.Range('E30').Value = '  -1.19%  '
$ws.Range('E31').Value = '  -4.81%  '
$ws.Range('E32').Value = '  +3.02%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').Value = '18.31'
$ws.Range('E34').Value = '  -0.61%  '
$ws.Range('E35').Value = '  -2.42%  '
$ws.Range('D36').Value = '4.04'
$ws.Range('E36').Value = '  -1.02%  '
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('D38').Value = '0.799'
$ws.Range('E38').Value = '  -2.12%  '
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('D40').Value = '280.66'
$ws.Range('E40').Value = '  +1.26%  '
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('E42').Value = '  -3.92%  '
$ws.Range('D43').Value = '129.55'
$ws.Range('E43').Value = '  +3.47%  '
$ws.Range('D44').Value = '10.91'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').Value = '0.599'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '0.0924'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('E47').Value = '  -2.24%  '
$ws.Range('D48').Value = '0.0217'
$ws.Range('E48').Value = '  -0.94%  '
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('D50').Value = '1.749.47'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').Value = '0.982'
$ws.Range('E51').Value = '  -0.34%  '
